$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Shift the two existing data rows (old row 4 -> row 7, old row 5 -> row 8)
#    down by inserting three blank rows at rows 3-5. Inserting whole rows
#    (rather than rebuilding them) preserves the row-level formatting
#    (s="2" customFormat="1") that the old row 4 already carried.
# ---------------------------------------------------------------------------
$ws.Range("A3:A5").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2) New row 3: "22-10-204" / AWS Certified Developer Associate Course
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("G7").Copy()
$ws.Range("G3").PasteSpecial(-4122)

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "22-10-204"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "AWS Certified Developer Associate Course"
$ws.Range("E3").Value = "Getting started, Intro, IAM & CLI"
$ws.Range("F3").Value = 0.083333333333333329
$ws.Range("G3").Value = 0.75

# ---------------------------------------------------------------------------
# 3) New row 5: 23-Oct-2024 / Learning Concetps / EC2,ELB,VPC,S3
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("G7").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("A5").Value = 2
$ws.Range("B5").Value2 = 45588
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "Learning Concetps"
$ws.Range("E5").Value = "EC2,ELB,VPC,S3"
$ws.Range("F5").Value = 0.39583333333333331
$ws.Range("G5").Value = 0.78472222222222221

# ---------------------------------------------------------------------------
# 4) Row 8 (old row 5, shifted down) gets a new G8 end-time value
# ---------------------------------------------------------------------------
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = 0.45833333333333331

# ---------------------------------------------------------------------------
# 5) New rows 9, 10, 11 appended after row 8
# ---------------------------------------------------------------------------
$ws.Range("F7").Copy()
$ws.Range("F9:F11").PasteSpecial(-4122)
$ws.Range("G7").Copy()
$ws.Range("G9:G11").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D9:D11").PasteSpecial(-4122)

# Rows 8, 9 and 10 keep empty (but styled) A/B cells since the Date / SL.NO
# are not repeated for continuation tasks on the same day.
$ws.Range("A8").Value = 0
$ws.Range("A8").ClearContents()
$ws.Range("B8").Value = 0
$ws.Range("B8").ClearContents()
$ws.Range("A9").Value = 0
$ws.Range("A9").ClearContents()
$ws.Range("B9").Value = 0
$ws.Range("B9").ClearContents()
$ws.Range("A10").Value = 0
$ws.Range("A10").ClearContents()
$ws.Range("B10").Value = 0
$ws.Range("B10").ClearContents()

# Row 9: GIT
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "GIT"
$ws.Range("E9").Value = "installing, created respository for my progress report and update my report"
$ws.Range("F9").Value = 0.47916666666666669
$ws.Range("G9").Value = 0.50694444444444442

# Row 10: Lambda Functions
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = "Lambda Functions"
$ws.Range("E10").Value = "Syncronous Invocations, services,lambda & ALB(HTTP-JSON, vice versa)"
$ws.Range("F10").Value = 0.50694444444444442
$ws.Range("G10").Value = 0.55208333333333337

# Row 11: Lambda Functions / S3 event notifications ... (no SL.NO / Date repeated)
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = "Lambda Functions"
$ws.Range("E11").Value = "S3 event notifications, event source mapping,strea,s, queues,error handling, mapper scalling"
$ws.Range("F11").Value = 0.59027777777777779
$ws.Range("G11").Value = 0.63680555555555551

# ---------------------------------------------------------------------------
# 6) Column width adjustments (best achievable match under the host's
#    character-width quantisation) and selection
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 9.25
$ws.Columns.Item(4).ColumnWidth = 30.92
$ws.Columns.Item(5).ColumnWidth = 62.08
$ws.Columns.Item(6).ColumnWidth = 11.08
$ws.Columns.Item(7).ColumnWidth = 9.92

$ws.Range("D16").Select()
